# Standardize "Multi KDMA" wording in the definition row (row 2) of Sheet1.
# The source spreadsheet used inconsistent casing/spacing ("multi-kdma",
# "multikdma") when referring to the Multi KDMA condition; this normalizes
# all occurrences to "Multi KDMA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("K2", "AA2", "AD2", "AG2", "AI2", "AL2", "AM2")

foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    $text = $text -replace "multi-kdma", "Multi KDMA"
    $text = $text -replace "multikdma", "Multi KDMA"
    $cell.Value2 = $text
}
